$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.222.55"
$ws.Range("E2").Value = "  -0.63%  "
Set-TextValue "D3" "1.784.52"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue "D5" "334.67"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("E6").Value = "  +0.08%  "
Set-TextValue "D7" "0.3784"
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("E8").Value = "  -3.02%  "
Set-TextValue "D9" "48.35"
$ws.Range("E9").Value = "  -4.44%  "
Set-TextValue "D10" "1.198"
$ws.Range("E10").Value = "  -3.78%  "
Set-TextValue "D11" "0.07499"
$ws.Range("E11").Value = "  -4.05%  "
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  -0.02%  "
Set-TextValue "D13" "21.76"
$ws.Range("E13").Value = "  -4.14%  "
Set-TextValue "D14" "6.477"
$ws.Range("E14").Value = "  -3.05%  "
Set-TextValue "D15" "1.790.71"
$ws.Range("E15").Value = "  -1.31%  "
Set-TextValue "D16" "7.105"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("E17").Value = "  -3.48%  "
Set-TextValue "D18" "0.06669"
Set-TextValue "D19" "83.80"
$ws.Range("E19").Value = "  -3.66%  "
Set-TextValue "D20" "1.001"
$ws.Range("E20").Value = "  +0.07%  "
Set-TextValue "D21" "6.622"
$ws.Range("E21").Value = "  +0.45%  "
Set-TextValue "D22" "17.35"
$ws.Range("E22").Value = "  -3.61%  "
Set-TextValue "D23" "27.224.38"
$ws.Range("E23").Value = "  -0.60%  "
Set-TextValue "D24" "12.40"
$ws.Range("E24").Value = "  -6.43%  "
$ws.Range("E25").Value = "  -2.04%  "
Set-TextValue "D26" "1.502"
$ws.Range("E26").Value = "  -1.25%  "
Set-TextValue "D27" "2.548"
$ws.Range("E27").Value = "  -7.81%  "
Set-TextValue "D28" "21.32"
$ws.Range("E28").Value = "  -3.52%  "
Set-TextValue "D29" "153.87"
$ws.Range("E29").Value = "  -0.21%  "
Set-TextValue "D30" "1.989.69"
$ws.Range("E30").Value = "  -1.37%  "
Set-TextValue "D31" "134.03"
$ws.Range("E31").Value = "  -2.73%  "
Set-TextValue "D32" "4.022"
$ws.Range("E32").Value = "  -2.62%  "
Set-TextValue "D33" "6.109"
$ws.Range("E33").Value = "  -5.55%  "
Set-TextValue "D34" "0.08696"
$ws.Range("E34").Value = "  -1.70%  "
Set-TextValue "D35" "13.29"
$ws.Range("E35").Value = "  -4.82%  "
Set-TextValue "D36" "1.662"
$ws.Range("E36").Value = "  -3.72%  "
Set-TextValue "D37" "0.6958"
$ws.Range("E37").Value = "  -3.58%  "
Set-TextValue "D38" "5.461"
$ws.Range("E38").Value = "  -3.95%  "
Set-TextValue "D39" "0.2204"
$ws.Range("E39").Value = "  -3.35%  "
Set-TextValue "D40" "8.826"
$ws.Range("E40").Value = "  -2.99%  "
Set-TextValue "D41" "0.06334"
$ws.Range("E41").Value = "  -3.99%  "
Set-TextValue "D42" "0.02338"
$ws.Range("E42").Value = "  -4.03%  "
Set-TextValue "D43" "1.243"
$ws.Range("E43").Value = "  -2.15%  "
Set-TextValue "D44" "14.41"
$ws.Range("E44").Value = "  -3.99%  "
Set-TextValue "D45" "0.6517"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("E47").Value = "  -3.14%  "
Set-TextValue "D48" "2.148"
$ws.Range("E48").Value = "  -2.70%  "
Set-TextValue "D49" "129.42"
$ws.Range("E49").Value = "  -3.51%  "
Set-TextValue "D50" "0.07141"
$ws.Range("E50").Value = "  -2.84%  "
Set-TextValue "D51" "79.30"
$ws.Range("E51").Value = "  -2.61%  "
